# MtomToCrss_Annual.xlsx — MRM re-run output update
# The "Trace32" sheet's F3 (and its dependent total, H3) were refreshed with
# the latest model-run numbers; every other sheet/cell is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trace32")

$ws.Range("F3").Value = 8855471.0505504292
$ws.Range("H3").Value = 9039420.770550428
